$d = $word.ActiveDocument

$d.Content.Find.Execute("695×4=2780", $true, $true, $false, $false, $false, $true, 1, $false, "781×6=4686", 2) | Out-Null
$d.Content.Find.Execute("455×4=1820", $true, $true, $false, $false, $false, $true, 1, $false, "214×8=1712", 2) | Out-Null
$d.Content.Find.Execute("429×7=3003", $true, $true, $false, $false, $false, $true, 1, $false, "842×5=4210", 2) | Out-Null
$d.Content.Find.Execute("325×2=650", $true, $true, $false, $false, $false, $true, 1, $false, "473×3=1419", 2) | Out-Null
$d.Content.Find.Execute("818×4=3272", $true, $true, $false, $false, $false, $true, 1, $false, "787×8=6296", 2) | Out-Null
$d.Content.Find.Execute("106×6=636", $true, $true, $false, $false, $false, $true, 1, $false, "916×9=8244", 2) | Out-Null
$d.Content.Find.Execute("150×3=450", $true, $true, $false, $false, $false, $true, 1, $false, "479×8=3832", 2) | Out-Null
$d.Content.Find.Execute("400×9=3600", $true, $true, $false, $false, $false, $true, 1, $false, "179×5=895", 2) | Out-Null
$d.Content.Find.Execute("849×8=6792", $true, $true, $false, $false, $false, $true, 1, $false, "141×7=987", 2) | Out-Null
$d.Content.Find.Execute("703×6=4218", $true, $true, $false, $false, $false, $true, 1, $false, "578×3=1734", 2) | Out-Null
$d.Content.Find.Execute("759×6=4554", $true, $true, $false, $false, $false, $true, 1, $false, "307×8=2456", 2) | Out-Null
$d.Content.Find.Execute("131×7=917", $true, $true, $false, $false, $false, $true, 1, $false, "512×6=3072", 2) | Out-Null
$d.Content.Find.Execute("493×3=1479", $true, $true, $false, $false, $false, $true, 1, $false, "499×4=1996", 2) | Out-Null
$d.Content.Find.Execute("380×9=3420", $true, $true, $false, $false, $false, $true, 1, $false, "870×4=3480", 2) | Out-Null
$d.Content.Find.Execute("992×4=3968", $true, $true, $false, $false, $false, $true, 1, $false, "903×2=1806", 2) | Out-Null
$d.Content.Find.Execute("735×3=2205", $true, $true, $false, $false, $false, $true, 1, $false, "756×2=1512", 2) | Out-Null
$d.Content.Find.Execute("297×4=1188", $true, $true, $false, $false, $false, $true, 1, $false, "749×7=5243", 2) | Out-Null
$d.Content.Find.Execute("659×2=1318", $true, $true, $false, $false, $false, $true, 1, $false, "675×4=2700", 2) | Out-Null
$d.Content.Find.Execute("986×2=1972", $true, $true, $false, $false, $false, $true, 1, $false, "610×4=2440", 2) | Out-Null
$d.Content.Find.Execute("724×7=5068", $true, $true, $false, $false, $false, $true, 1, $false, "977×7=6839", 2) | Out-Null
$d.Content.Find.Execute("212×7=1484", $true, $true, $false, $false, $false, $true, 1, $false, "438×2=876", 2) | Out-Null
$d.Content.Find.Execute("316×4=1264", $true, $true, $false, $false, $false, $true, 1, $false, "785×2=1570", 2) | Out-Null
$d.Content.Find.Execute("485×5=2425", $true, $true, $false, $false, $false, $true, 1, $false, "418×5=2090", 2) | Out-Null
$d.Content.Find.Execute("560×3=1680", $true, $true, $false, $false, $false, $true, 1, $false, "714×6=4284", 2) | Out-Null
$d.Content.Find.Execute("213×3=639", $true, $true, $false, $false, $false, $true, 1, $false, "361×7=2527", 2) | Out-Null

Write-Output "Done replacing values"
